# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (F column) figures and one cover image link on the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$sheet1Updates = @{
    "F2"  = 1277
    "F3"  = 1168
    "F5"  = 105
    "F7"  = 647
    "F8"  = 93
    "F11" = 2314
    "F12" = 1582
    "F13" = 1321
    "F15" = 230
    "F16" = 536
    "F17" = 746
    "F18" = 35
    "F19" = 281
    "F20" = 1085
    "F22" = 12
    "F24" = 4541
    "F25" = 209
    "F26" = 30
    "F28" = 124
    "F29" = 199
    "F31" = 8
    "F32" = 656
    "F36" = 235
    "F37" = 367
    "F38" = 955
    "F39" = 125
    "F40" = 91
    "F41" = 136
}

foreach ($addr in $sheet1Updates.Keys) {
    $ws1.Range($addr).Value = $sheet1Updates[$addr]
}

$ws1.Range("I32").Value = "//i0.hdslb.com/bfs/openplatform/202402/RDI807mS1708410823039.jpeg"

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Updates = @{
    "F2"  = 1277
    "F5"  = 1168
    "F9"  = 105
    "F11" = 647
    "F12" = 93
    "F18" = 2314
    "F19" = 1582
    "F20" = 1321
    "F22" = 230
    "F23" = 536
    "F25" = 746
    "F26" = 35
    "F27" = 281
    "F28" = 1085
    "F29" = 12
    "F30" = 4541
    "F31" = 209
    "F32" = 30
    "F34" = 124
    "F35" = 199
    "F37" = 8
    "F38" = 656
    "F41" = 367
    "F42" = 955
    "F43" = 125
    "F44" = 91
    "F45" = 136
}

foreach ($addr in $sheet4Updates.Keys) {
    $ws4.Range($addr).Value = $sheet4Updates[$addr]
}

$ws4.Range("I38").Value = "//i0.hdslb.com/bfs/openplatform/202402/RDI807mS1708410823039.jpeg"
